$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 349.35715
$ws.Range("J32").Value = 349.0909
$ws.Range("L32").Value = 349.0909
$ws.Range("N32").Value = -1001.0909

$ws.Range("H40").Value = 2753.1516
$ws.Range("J40").Value = 3386
$ws.Range("L40").Value = 3386
$ws.Range("N40").Value = -3736

$ws.Range("H129").Value = 2584.6833
$ws.Range("I129").Value = 8270.77
$ws.Range("J129").Value = 1011.93616
$ws.Range("K129").Value = 24812.31
$ws.Range("L129").Value = 3035.80848
$ws.Range("M129").Value = -19812.31
$ws.Range("N129").Value = -13035.80848

$ws.Range("H135").Value = 1327.3667
$ws.Range("I135").Value = 622.95
$ws.Range("J135").Value = 2736.2
$ws.Range("K135").Value = 5606.55
$ws.Range("L135").Value = 24625.8
$ws.Range("M135").Value = -3071.55
$ws.Range("N135").Value = -29695.8

$ws.Range("H137").Value = 1579.619
$ws.Range("I137").Value = 1719.3334
$ws.Range("J137").Value = 1230.3334
$ws.Range("K137").Value = 5158.0002
$ws.Range("L137").Value = 3691.0002
$ws.Range("M137").Value = -2608.0002
$ws.Range("N137").Value = -8791.0002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 26915.783
$ws.Range("I32").Value = 5893.4263
$ws.Range("J32").Value = 265169.16
$ws.Range("K32").Value = 5893.4263
$ws.Range("L32").Value = 265169.16
$ws.Range("M32").Value = -5606.4263
$ws.Range("N32").Value = -265743.16

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 21860.4
$ws.Range("J35").Value = 21860.4
$ws.Range("L35").Value = 21860.4
$ws.Range("N35").Value = -22480.4

$ws.Range("H105").Value = 61006.47
$ws.Range("I105").Value = 73743.57000000001
$ws.Range("J105").Value = 1566.6666
$ws.Range("K105").Value = 73743.57000000001
$ws.Range("L105").Value = 1566.6666
$ws.Range("M105").Value = -71996.57000000001
$ws.Range("N105").Value = -5060.6666

$ws.Range("H134").Value = 6342.6
$ws.Range("I134").Value = 6502.4
$ws.Range("J134").Value = 6182.8
$ws.Range("K134").Value = 19507.2
$ws.Range("L134").Value = 18548.4
$ws.Range("M134").Value = -16972.2
$ws.Range("N134").Value = -23618.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1138.6471
$ws.Range("I5").Value = 821.6
$ws.Range("J5").Value = 1388.9474
$ws.Range("K5").Value = 2464.8
$ws.Range("L5").Value = 4166.8422
$ws.Range("M5").Value = -2352.8
$ws.Range("N5").Value = -4390.8422

$ws.Range("H23").Value = 1106.7368
$ws.Range("I23").Value = 2232
$ws.Range("J23").Value = 704.8570999999999
$ws.Range("K23").Value = 6696
$ws.Range("L23").Value = 2114.5713
$ws.Range("M23").Value = -6461
$ws.Range("N23").Value = -2584.5713

$ws.Range("H38").Value = 34.92857
$ws.Range("J38").Value = 58.8
$ws.Range("L38").Value = 176.4
$ws.Range("N38").Value = -870.4

$ws.Range("H58").Value = 1900
$ws.Range("J58").Value = 2900
$ws.Range("L58").Value = 8700
$ws.Range("N58").Value = -8956

$ws.Range("H107").Value = 535189.4399999999
$ws.Range("I107").Value = 571.38464
$ws.Range("K107").Value = 1714.15392
$ws.Range("M107").Value = 205.84608

$ws.Range("H113").Value = 667.61536
$ws.Range("I113").Value = 675.8
$ws.Range("J113").Value = 662.5
$ws.Range("K113").Value = 2027.4
$ws.Range("L113").Value = 1987.5
$ws.Range("M113").Value = 142.6000000000001
$ws.Range("N113").Value = -6327.5

$ws.Range("H131").Value = 796.85
$ws.Range("I131").Value = 377.66666
$ws.Range("J131").Value = 809.81445
$ws.Range("K131").Value = 1132.99998
$ws.Range("L131").Value = 2429.44335
$ws.Range("M131").Value = 3907.00002
$ws.Range("N131").Value = -12509.44335

$ws.Range("H135").Value = 1138.6471
$ws.Range("I135").Value = 821.6
$ws.Range("J135").Value = 1388.9474
$ws.Range("K135").Value = 7394.400000000001
$ws.Range("L135").Value = 12500.5266
$ws.Range("M135").Value = -4859.400000000001
$ws.Range("N135").Value = -17570.5266

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 34579.547
$ws.Range("I70").Value = 42628.383
$ws.Range("J70").Value = 4683.857
$ws.Range("K70").Value = 42628.383
$ws.Range("L70").Value = 4683.857
$ws.Range("M70").Value = -42358.383
$ws.Range("N70").Value = -5223.857

$ws.Range("H73").Value = 34579.547
$ws.Range("I73").Value = 42628.383
$ws.Range("J73").Value = 4683.857
$ws.Range("K73").Value = 42628.383
$ws.Range("L73").Value = 4683.857
$ws.Range("M73").Value = -41692.383
$ws.Range("N73").Value = -6555.857

$ws.Range("H126").Value = 3310.3572
$ws.Range("I126").Value = 3149.6155
$ws.Range("K126").Value = 9448.8465
$ws.Range("M126").Value = -6978.8465

$ws.Range("H132").Value = 3139.45
$ws.Range("I132").Value = 2996.5715
$ws.Range("J132").Value = 3472.8333
$ws.Range("K132").Value = 8989.7145
$ws.Range("L132").Value = 10418.4999
$ws.Range("M132").Value = -6459.7145
$ws.Range("N132").Value = -15478.4999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 375375
$ws.Range("I2").Value = 375375
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 375375
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -375263
$ws.Range("N2").ClearContents()

$ws.Range("H61").Value = 1709.0869
$ws.Range("J61").Value = 1900.3636
$ws.Range("L61").Value = 1900.3636
$ws.Range("N61").Value = -2304.3636

$ws.Range("H113").Value = 1709.0869
$ws.Range("J113").Value = 1900.3636
$ws.Range("L113").Value = 1900.3636
$ws.Range("N113").Value = -6240.3636

$ws.Range("H134").Value = 60704.145
$ws.Range("J134").Value = 60704.145
$ws.Range("L134").Value = 60704.145
$ws.Range("N134").Value = -70844.14499999999

$ws.Range("H136").Value = 2429
$ws.Range("I136").Value = 2055.7144
$ws.Range("J136").Value = 3300
$ws.Range("K136").Value = 6167.1432
$ws.Range("L136").Value = 9900
$ws.Range("M136").Value = -3617.1432
$ws.Range("N136").Value = -15000

$ws.Range("H138").Value = 64001.8
$ws.Range("J138").Value = 64001.8
$ws.Range("L138").Value = 64001.8
$ws.Range("N138").Value = -74281.8

$ws.Range("H140").Value = 88405.8
$ws.Range("J140").Value = 88405.8
$ws.Range("L140").Value = 88405.8
$ws.Range("N140").Value = -98765.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H56").Value = 31451.715
$ws.Range("J56").Value = 36027
$ws.Range("L56").Value = 36027
$ws.Range("N56").Value = -37455

$ws.Range("H132").Value = 4402
$ws.Range("I132").Value = 4513.1113
$ws.Range("J132").Value = 3402
$ws.Range("K132").Value = 13539.3339
$ws.Range("L132").Value = 10206
$ws.Range("M132").Value = -11009.3339
$ws.Range("N132").Value = -15266

$ws.Range("H136").Value = 2318
$ws.Range("I136").Value = 668.63635
$ws.Range("J136").Value = 4333.8887
$ws.Range("K136").Value = 2005.90905
$ws.Range("L136").Value = 13001.6661
$ws.Range("M136").Value = 544.09095
$ws.Range("N136").Value = -18101.6661

$ws.Range("H138").Value = 66180
$ws.Range("J138").Value = 66180
$ws.Range("L138").Value = 66180
$ws.Range("N138").Value = -76460
